$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# "ODI Batting Extra" gained a full extra scrape of 34 additional
# MATCH_CODE rows. Insert 34 blank rows above the existing data (this
# shifts the previously-first data row, 4166 @ row 2, down to row 36,
# preserving all of its original cell content/formatting untouched)
# and then populate the newly inserted rows 2-35 with the new data.
$ws.Rows.Item(2).Resize(34).Insert()

# Row Insert copies the formatting of the row above (the bold header
# row here) into the freshly inserted rows; strip that back off so the
# new data rows stay plain, like all the other body rows.
$ws.Rows.Item(2).Resize(34).ClearFormats()

# Row 2 (MATCH_CODE 3898)
$ws.Cells.Item(2,1).Value = "'3898"
$ws.Cells.Item(2,2).Value = 6
$ws.Cells.Item(2,3).Value = "'0"
$ws.Cells.Item(2,4).Value = "'0"
$ws.Cells.Item(2,5).Value = "'0.38%"
$ws.Cells.Item(2,6).Value = "'NO"

# Row 3 (MATCH_CODE 3923)
$ws.Cells.Item(3,1).Value = "'3923"
$ws.Cells.Item(3,2).Value = ""
$ws.Cells.Item(3,3).Value = ""
$ws.Cells.Item(3,4).Value = ""
$ws.Cells.Item(3,5).Value = ""
$ws.Cells.Item(3,6).Value = "'NO"

# Row 4 (MATCH_CODE 3924)
$ws.Cells.Item(4,1).Value = "'3924"
$ws.Cells.Item(4,2).Value = 7
$ws.Cells.Item(4,3).Value = "'3"
$ws.Cells.Item(4,4).Value = "'0"
$ws.Cells.Item(4,5).Value = "'15.05%"
$ws.Cells.Item(4,6).Value = "'NO"

# Row 5 (MATCH_CODE 3927)
$ws.Cells.Item(5,1).Value = "'3927"
$ws.Cells.Item(5,2).Value = ""
$ws.Cells.Item(5,3).Value = ""
$ws.Cells.Item(5,4).Value = ""
$ws.Cells.Item(5,5).Value = ""
$ws.Cells.Item(5,6).Value = "'NO"

# Row 6 (MATCH_CODE 3929)
$ws.Cells.Item(6,1).Value = "'3929"
$ws.Cells.Item(6,2).Value = ""
$ws.Cells.Item(6,3).Value = ""
$ws.Cells.Item(6,4).Value = ""
$ws.Cells.Item(6,5).Value = ""
$ws.Cells.Item(6,6).Value = "'NO"

# Row 7 (MATCH_CODE 3931)
$ws.Cells.Item(7,1).Value = "'3931"
$ws.Cells.Item(7,2).Value = 5
$ws.Cells.Item(7,3).Value = "'1"
$ws.Cells.Item(7,4).Value = "'0"
$ws.Cells.Item(7,5).Value = "'6.53%"
$ws.Cells.Item(7,6).Value = "'NO"

# Row 8 (MATCH_CODE 3937)
$ws.Cells.Item(8,1).Value = "'3937"
$ws.Cells.Item(8,2).Value = 5
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(8,4).Value = ""
$ws.Cells.Item(8,5).Value = ""
$ws.Cells.Item(8,6).Value = "'NO"

# Row 9 (MATCH_CODE 3940)
$ws.Cells.Item(9,1).Value = "'3940"
$ws.Cells.Item(9,2).Value = ""
$ws.Cells.Item(9,3).Value = ""
$ws.Cells.Item(9,4).Value = ""
$ws.Cells.Item(9,5).Value = ""
$ws.Cells.Item(9,6).Value = "'NO"

# Row 10 (MATCH_CODE 3942)
$ws.Cells.Item(10,1).Value = "'3942"
$ws.Cells.Item(10,2).Value = 6
$ws.Cells.Item(10,3).Value = "'5"
$ws.Cells.Item(10,4).Value = "'2"
$ws.Cells.Item(10,5).Value = "'23.29%"
$ws.Cells.Item(10,6).Value = "'NO"

# Row 11 (MATCH_CODE 3945)
$ws.Cells.Item(11,1).Value = "'3945"
$ws.Cells.Item(11,2).Value = 6
$ws.Cells.Item(11,3).Value = "'4"
$ws.Cells.Item(11,4).Value = "'1"
$ws.Cells.Item(11,5).Value = "'9.43%"
$ws.Cells.Item(11,6).Value = "'NO"

# Row 12 (MATCH_CODE 3947)
$ws.Cells.Item(12,1).Value = "'3947"
$ws.Cells.Item(12,2).Value = ""
$ws.Cells.Item(12,3).Value = ""
$ws.Cells.Item(12,4).Value = ""
$ws.Cells.Item(12,5).Value = ""
$ws.Cells.Item(12,6).Value = "'NO"

# Row 13 (MATCH_CODE 3950)
$ws.Cells.Item(13,1).Value = "'3950"
$ws.Cells.Item(13,2).Value = ""
$ws.Cells.Item(13,3).Value = ""
$ws.Cells.Item(13,4).Value = ""
$ws.Cells.Item(13,5).Value = ""
$ws.Cells.Item(13,6).Value = "'NO"

# Row 14 (MATCH_CODE 3966)
$ws.Cells.Item(14,1).Value = "'3966"
$ws.Cells.Item(14,2).Value = ""
$ws.Cells.Item(14,3).Value = ""
$ws.Cells.Item(14,4).Value = ""
$ws.Cells.Item(14,5).Value = ""
$ws.Cells.Item(14,6).Value = "'NO"

# Row 15 (MATCH_CODE 3967)
$ws.Cells.Item(15,1).Value = "'3967"
$ws.Cells.Item(15,2).Value = 4
$ws.Cells.Item(15,3).Value = "'6"
$ws.Cells.Item(15,4).Value = "'2"
$ws.Cells.Item(15,5).Value = "'15.08%"
$ws.Cells.Item(15,6).Value = "'NO"

# Row 16 (MATCH_CODE 3968)
$ws.Cells.Item(16,1).Value = "'3968"
$ws.Cells.Item(16,2).Value = ""
$ws.Cells.Item(16,3).Value = ""
$ws.Cells.Item(16,4).Value = ""
$ws.Cells.Item(16,5).Value = ""
$ws.Cells.Item(16,6).Value = "'NO"

# Row 17 (MATCH_CODE 3972)
$ws.Cells.Item(17,1).Value = "'3972"
$ws.Cells.Item(17,2).Value = 2
$ws.Cells.Item(17,3).Value = "'5"
$ws.Cells.Item(17,4).Value = "'0"
$ws.Cells.Item(17,5).Value = "'14.55%"
$ws.Cells.Item(17,6).Value = "'NO"

# Row 18 (MATCH_CODE 3973)
$ws.Cells.Item(18,1).Value = "'3973"
$ws.Cells.Item(18,2).Value = ""
$ws.Cells.Item(18,3).Value = ""
$ws.Cells.Item(18,4).Value = ""
$ws.Cells.Item(18,5).Value = ""
$ws.Cells.Item(18,6).Value = "'NO"

# Row 19 (MATCH_CODE 3975)
$ws.Cells.Item(19,1).Value = "'3975"
$ws.Cells.Item(19,2).Value = 5
$ws.Cells.Item(19,3).Value = "'4"
$ws.Cells.Item(19,4).Value = "'0"
$ws.Cells.Item(19,5).Value = "'8.68%"
$ws.Cells.Item(19,6).Value = "'NO"

# Row 20 (MATCH_CODE 3977)
$ws.Cells.Item(20,1).Value = "'3977"
$ws.Cells.Item(20,2).Value = 4
$ws.Cells.Item(20,3).Value = "'2"
$ws.Cells.Item(20,4).Value = "'4"
$ws.Cells.Item(20,5).Value = "'14.45%"
$ws.Cells.Item(20,6).Value = "'NO"

# Row 21 (MATCH_CODE 3981)
$ws.Cells.Item(21,1).Value = "'3981"
$ws.Cells.Item(21,2).Value = 2
$ws.Cells.Item(21,3).Value = "'9"
$ws.Cells.Item(21,4).Value = "'3"
$ws.Cells.Item(21,5).Value = "'34.69%"
$ws.Cells.Item(21,6).Value = "'NO"

# Row 22 (MATCH_CODE 3984)
$ws.Cells.Item(22,1).Value = "'3984"
$ws.Cells.Item(22,2).Value = 2
$ws.Cells.Item(22,3).Value = "'1"
$ws.Cells.Item(22,4).Value = "'0"
$ws.Cells.Item(22,5).Value = "'1.79%"
$ws.Cells.Item(22,6).Value = "'NO"

# Row 23 (MATCH_CODE 3988)
$ws.Cells.Item(23,1).Value = "'3988"
$ws.Cells.Item(23,2).Value = 4
$ws.Cells.Item(23,3).Value = "'5"
$ws.Cells.Item(23,4).Value = "'0"
$ws.Cells.Item(23,5).Value = "'20.62%"
$ws.Cells.Item(23,6).Value = "'NO"

# Row 24 (MATCH_CODE 4032)
$ws.Cells.Item(24,1).Value = "'4032"
$ws.Cells.Item(24,2).Value = 6
$ws.Cells.Item(24,3).Value = ""
$ws.Cells.Item(24,4).Value = ""
$ws.Cells.Item(24,5).Value = ""
$ws.Cells.Item(24,6).Value = "'NO"

# Row 25 (MATCH_CODE 4035)
$ws.Cells.Item(25,1).Value = "'4035"
$ws.Cells.Item(25,2).Value = 6
$ws.Cells.Item(25,3).Value = ""
$ws.Cells.Item(25,4).Value = ""
$ws.Cells.Item(25,5).Value = ""
$ws.Cells.Item(25,6).Value = "'NO"

# Row 26 (MATCH_CODE 4041)
$ws.Cells.Item(26,1).Value = "'4041"
$ws.Cells.Item(26,2).Value = ""
$ws.Cells.Item(26,3).Value = ""
$ws.Cells.Item(26,4).Value = ""
$ws.Cells.Item(26,5).Value = ""
$ws.Cells.Item(26,6).Value = "'NO"

# Row 27 (MATCH_CODE 4067)
$ws.Cells.Item(27,1).Value = "'4067"
$ws.Cells.Item(27,2).Value = 4
$ws.Cells.Item(27,3).Value = "'0"
$ws.Cells.Item(27,4).Value = "'0"
$ws.Cells.Item(27,5).Value = "'3.65%"
$ws.Cells.Item(27,6).Value = "'NO"

# Row 28 (MATCH_CODE 4069)
$ws.Cells.Item(28,1).Value = "'4069"
$ws.Cells.Item(28,2).Value = 4
$ws.Cells.Item(28,3).Value = "'5"
$ws.Cells.Item(28,4).Value = "'0"
$ws.Cells.Item(28,5).Value = "'19.31%"
$ws.Cells.Item(28,6).Value = "'NO"

# Row 29 (MATCH_CODE 4071)
$ws.Cells.Item(29,1).Value = "'4071"
$ws.Cells.Item(29,2).Value = 5
$ws.Cells.Item(29,3).Value = "'0"
$ws.Cells.Item(29,4).Value = "'0"
$ws.Cells.Item(29,5).Value = "'1.37%"
$ws.Cells.Item(29,6).Value = "'NO"

# Row 30 (MATCH_CODE 4074)
$ws.Cells.Item(30,1).Value = "'4074"
$ws.Cells.Item(30,2).Value = ""
$ws.Cells.Item(30,3).Value = ""
$ws.Cells.Item(30,4).Value = ""
$ws.Cells.Item(30,5).Value = ""
$ws.Cells.Item(30,6).Value = "'NO"

# Row 31 (MATCH_CODE 4076)
$ws.Cells.Item(31,1).Value = "'4076"
$ws.Cells.Item(31,2).Value = 5
$ws.Cells.Item(31,3).Value = "'4"
$ws.Cells.Item(31,4).Value = "'0"
$ws.Cells.Item(31,5).Value = "'17.36%"
$ws.Cells.Item(31,6).Value = "'NO"

# Row 32 (MATCH_CODE 4108)
$ws.Cells.Item(32,1).Value = "'4108"
$ws.Cells.Item(32,2).Value = ""
$ws.Cells.Item(32,3).Value = ""
$ws.Cells.Item(32,4).Value = ""
$ws.Cells.Item(32,5).Value = ""
$ws.Cells.Item(32,6).Value = "'NO"

# Row 33 (MATCH_CODE 4115)
$ws.Cells.Item(33,1).Value = "'4115"
$ws.Cells.Item(33,2).Value = 4
$ws.Cells.Item(33,3).Value = "'0"
$ws.Cells.Item(33,4).Value = "'0"
$ws.Cells.Item(33,5).Value = "'2.59%"
$ws.Cells.Item(33,6).Value = "'NO"

# Row 34 (MATCH_CODE 4123)
$ws.Cells.Item(34,1).Value = "'4123"
$ws.Cells.Item(34,2).Value = 2
$ws.Cells.Item(34,3).Value = "'15"
$ws.Cells.Item(34,4).Value = "'0"
$ws.Cells.Item(34,5).Value = "'48.73%"
$ws.Cells.Item(34,6).Value = "'NO"

# Row 35 (MATCH_CODE 4125)
$ws.Cells.Item(35,1).Value = "'4125"
$ws.Cells.Item(35,2).Value = 2
$ws.Cells.Item(35,3).Value = "'2"
$ws.Cells.Item(35,4).Value = "'0"
$ws.Cells.Item(35,5).Value = "'8.91%"
$ws.Cells.Item(35,6).Value = "'NO"
